$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of portfolio data for 2025-08-30 (same figures as the prior day).
# Force the Date column to be treated as literal text (matching the rest of
# column A, which stores dates as plain strings rather than date serials),
# then restore the default "Normal" style so no extra formatting is
# introduced on the new cell.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2025-08-30"
$ws.Range("A15").Style = "Normal"

$ws.Range("B15").Value = 56.43000030517578
$ws.Range("C15").Value = 669
$ws.Range("D15").Value = 313.9500122070312
